# Actualizacion a mejor modelo
# Update predicted 'sdg' category values (column B) for specific rows
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    81 = 3
    159 = 5
    188 = 5
    287 = 5
    294 = 4
    316 = 5
    363 = 5
    368 = 5
    402 = 4
    407 = 5
    413 = 3
    427 = 4
    463 = 3
    475 = 5
    479 = 4
    503 = 4
    524 = 4
    552 = 5
    572 = 4
    625 = 4
    644 = 4
    657 = 5
    690 = 5
    728 = 5
    863 = 5
    868 = 5
    897 = 5
    904 = 4
    986 = 3
    1017 = 5
    1040 = 3
    1046 = 3
    1096 = 5
    1137 = 3
    1183 = 3
    1184 = 4
    1189 = 5
    1191 = 4
    1209 = 3
    1213 = 4
    1255 = 4
    1311 = 3
    1364 = 3
    1367 = 5
    1372 = 3
    1383 = 4
    1404 = 5
    1416 = 4
    1420 = 5
    1495 = 4
    1513 = 5
    1523 = 3
    1535 = 4
    1580 = 4
    1649 = 4
    1675 = 3
    1698 = 5
    1731 = 5
    1733 = 5
    1754 = 4
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item([int]$row, 2).Value = $updates[$row]
}

